$wb = $excel.ActiveWorkbook

# ALC row 15: Morning Glass of Ether
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 11906396
$ws.Cells.Item(15, 9).Value = 11906396
$ws.Cells.Item(15, 11).Value = 35719188
$ws.Cells.Item(15, 13).Value = -35719019

# ALC row 138: All-night Crafting
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1452.1464
$ws.Cells.Item(138, 9).Value = 620.0278
$ws.Cells.Item(138, 10).Value = 7443.4
$ws.Cells.Item(138, 11).Value = 1860.0834
$ws.Cells.Item(138, 12).Value = 22330.2
$ws.Cells.Item(138, 13).Value = 3279.9166
$ws.Cells.Item(138, 14).Value = -32610.2

# ALC row 141: Remedy for Reason
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 2552.375
$ws.Cells.Item(141, 9).Value = 1124.9166
$ws.Cells.Item(141, 10).Value = 9689.666999999999
$ws.Cells.Item(141, 11).Value = 3374.7498
$ws.Cells.Item(141, 12).Value = 29069.001
$ws.Cells.Item(141, 13).Value = 1805.2502
$ws.Cells.Item(141, 14).Value = -39429.001

# ARM row 2: Ain't Got No Ingots
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 849.6667
$ws.Cells.Item(2, 9).Value = 849.6667
$ws.Cells.Item(2, 11).Value = 849.6667
$ws.Cells.Item(2, 13).Value = -736.6667

# ARM row 32: Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11467.043
$ws.Cells.Item(32, 9).Value = 2951.797
$ws.Cells.Item(32, 10).Value = 34969.12
$ws.Cells.Item(32, 11).Value = 2951.797
$ws.Cells.Item(32, 12).Value = 34969.12
$ws.Cells.Item(32, 13).Value = -2664.797
$ws.Cells.Item(32, 14).Value = -35543.12

# ARM row 116: No Scope
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 849.6667
$ws.Cells.Item(116, 9).Value = 849.6667
$ws.Cells.Item(116, 11).Value = 849.6667
$ws.Cells.Item(116, 13).Value = 1444.3333

# ARM row 123: The Armoire Is Open
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(123, 8).Value = 54300
$ws.Cells.Item(123, 10).Value = 54300
$ws.Cells.Item(123, 12).Value = 54300
$ws.Cells.Item(123, 14).Value = -64100

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1279.909
$ws.Cells.Item(132, 9).Value = 1181.6111
$ws.Cells.Item(132, 10).Value = 1722.25
$ws.Cells.Item(132, 11).Value = 3544.8333
$ws.Cells.Item(132, 12).Value = 5166.75
$ws.Cells.Item(132, 13).Value = -1014.8333
$ws.Cells.Item(132, 14).Value = -10226.75

# BSM row 3: Hells Bells
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 849.6667
$ws.Cells.Item(3, 9).Value = 849.6667
$ws.Cells.Item(3, 11).Value = 849.6667
$ws.Cells.Item(3, 13).Value = -735.6667

# BSM row 107: The Gold Experience
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1126.7241
$ws.Cells.Item(107, 9).Value = 892.0625
$ws.Cells.Item(107, 11).Value = 892.0625
$ws.Cells.Item(107, 13).Value = 1027.9375

# CRP row 31: Wall Not Found
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3046533
$ws.Cells.Item(31, 9).Value = 3466308.5
$ws.Cells.Item(31, 10).Value = 3161.125
$ws.Cells.Item(31, 11).Value = 3466308.5
$ws.Cells.Item(31, 12).Value = 3161.125
$ws.Cells.Item(31, 13).Value = -3466013.5
$ws.Cells.Item(31, 14).Value = -3751.125

# CRP row 34: Armoires of the Rich and Famous
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3046533
$ws.Cells.Item(34, 9).Value = 3466308.5
$ws.Cells.Item(34, 10).Value = 3161.125
$ws.Cells.Item(34, 11).Value = 3466308.5
$ws.Cells.Item(34, 12).Value = 3161.125
$ws.Cells.Item(34, 13).Value = -3466106.5
$ws.Cells.Item(34, 14).Value = -3565.125

# CRP row 86: Birch, Please
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 837334.5
$ws.Cells.Item(86, 9).Value = 1003401.4
$ws.Cells.Item(86, 10).Value = 7000
$ws.Cells.Item(86, 11).Value = 1003401.4
$ws.Cells.Item(86, 12).Value = 7000
$ws.Cells.Item(86, 13).Value = -1002278.4
$ws.Cells.Item(86, 14).Value = -9246

# CRP row 89: Built This City on Blocks and Soul (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 837334.5
$ws.Cells.Item(89, 9).Value = 1003401.4
$ws.Cells.Item(89, 10).Value = 7000
$ws.Cells.Item(89, 11).Value = 5017007
$ws.Cells.Item(89, 12).Value = 35000
$ws.Cells.Item(89, 13).Value = -5011391
$ws.Cells.Item(89, 14).Value = -46232

# CRP row 99: O Pine
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 3249757
$ws.Cells.Item(99, 9).Value = 5954594.5
$ws.Cells.Item(99, 10).Value = 3952
$ws.Cells.Item(99, 11).Value = 5954594.5
$ws.Cells.Item(99, 12).Value = 3952
$ws.Cells.Item(99, 13).Value = -5953096.5
$ws.Cells.Item(99, 14).Value = -6948

# CRP row 126: A Better Conductor
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 3249757
$ws.Cells.Item(126, 9).Value = 5954594.5
$ws.Cells.Item(126, 10).Value = 3952
$ws.Cells.Item(126, 11).Value = 17863783.5
$ws.Cells.Item(126, 12).Value = 11856
$ws.Cells.Item(126, 13).Value = -17861313.5
$ws.Cells.Item(126, 14).Value = -16796

# CUL row 8: Whip It
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 321.93332
$ws.Cells.Item(8, 9).Value = 321.93332
$ws.Cells.Item(8, 11).Value = 965.7999599999999
$ws.Cells.Item(8, 13).Value = -826.7999599999999

# CUL row 131: The Mountain Steeped
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 999.37836
$ws.Cells.Item(131, 10).Value = 1214.2593
$ws.Cells.Item(131, 12).Value = 3642.7779
$ws.Cells.Item(131, 14).Value = -13722.7779

# GSM row 113: Copious Crystal Cannons
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3416.6667
$ws.Cells.Item(113, 9).Value = 1375
$ws.Cells.Item(113, 10).Value = 7500
$ws.Cells.Item(113, 11).Value = 1375
$ws.Cells.Item(113, 12).Value = 7500
$ws.Cells.Item(113, 13).Value = 795
$ws.Cells.Item(113, 14).Value = -11840

# GSM row 122: Awarding Academic Excellence
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 4612.732
$ws.Cells.Item(122, 9).Value = 4165.0645
$ws.Cells.Item(122, 10).Value = 6000.5
$ws.Cells.Item(122, 11).Value = 12495.1935
$ws.Cells.Item(122, 12).Value = 18001.5
$ws.Cells.Item(122, 13).Value = -10045.1935
$ws.Cells.Item(122, 14).Value = -22901.5

# LTW row 7: Tan Before the Ban
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1948.75
$ws.Cells.Item(7, 9).Value = 1147.1666
$ws.Cells.Item(7, 10).Value = 2750.3333
$ws.Cells.Item(7, 11).Value = 1147.1666
$ws.Cells.Item(7, 12).Value = 2750.3333
$ws.Cells.Item(7, 13).Value = -1035.1666
$ws.Cells.Item(7, 14).Value = -2974.3333

# LTW row 122: Hell on Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 5556.75
$ws.Cells.Item(122, 9).Value = 5614.7144
$ws.Cells.Item(122, 11).Value = 16844.1432
$ws.Cells.Item(122, 13).Value = -14394.1432

# LTW row 126: Battered Books
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 1948.75
$ws.Cells.Item(126, 9).Value = 1147.1666
$ws.Cells.Item(126, 10).Value = 2750.3333
$ws.Cells.Item(126, 11).Value = 3441.4998
$ws.Cells.Item(126, 12).Value = 8250.999899999999
$ws.Cells.Item(126, 13).Value = -971.4998000000001
$ws.Cells.Item(126, 14).Value = -13190.9999

# WVR row 113: A Tender Table
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 423.3
$ws.Cells.Item(113, 9).Value = 474.125
$ws.Cells.Item(113, 10).Value = 220
$ws.Cells.Item(113, 11).Value = 1422.375
$ws.Cells.Item(113, 12).Value = 660
$ws.Cells.Item(113, 13).Value = 747.625
$ws.Cells.Item(113, 14).Value = -5000

# WVR row 115: Gloves Come in Handy
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(115, 8).Value = 20950
$ws.Cells.Item(115, 10).Value = 20950
$ws.Cells.Item(115, 12).Value = 20950
$ws.Cells.Item(115, 14).Value = -24084

# WVR row 122: Heavy Armoire
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1477.48
$ws.Cells.Item(122, 9).Value = 1109.3125
$ws.Cells.Item(122, 10).Value = 2132
$ws.Cells.Item(122, 11).Value = 3327.9375
$ws.Cells.Item(122, 12).Value = 6396
$ws.Cells.Item(122, 13).Value = -877.9375
$ws.Cells.Item(122, 14).Value = -11296

# WVR row 136: Weaving the Envelope
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1228.9608
$ws.Cells.Item(136, 9).Value = 818.0270400000001
$ws.Cells.Item(136, 10).Value = 2315
$ws.Cells.Item(136, 11).Value = 2454.08112
$ws.Cells.Item(136, 12).Value = 6945
$ws.Cells.Item(136, 13).Value = 95.91887999999972
$ws.Cells.Item(136, 14).Value = -12045

Write-Host "Applied scheduled runner updates to Zeromus_Profits sheets."
